# "object and arrays übungen"
#
# 1) Move the "Übung zu Array Methoden 'Rangliste'" slide (currently slide 13)
#    down to slide position 19 (right after the "selbe Tabelle" exercise slide).
# 2) Fix a typo on the "Arrays & Objekte" exercises slide: "die selbe" -> "dieselbe".
# 3) Refresh the cached datetimeFigureOut text (master + every layout) to the
#    current save date.

$p = $ppt.ActivePresentation

# --- 1. Reorder slides -----------------------------------------------------
$moved = $p.Slides.Item(13)
$moved.MoveTo(19)

# --- 2. Typo fix on the slide that is now at position 18 -------------------
$sl = $p.Slides.Item(18)
$shp = $sl.Shapes.Item(2)
$tr = $shp.TextFrame.TextRange
$full = $tr.Text
$needle = "4. Lege die selbe Tabelle als "
$idx = $full.IndexOf($needle)
$sub = $tr.Characters($idx + 1, $needle.Length)
$sub.Text = "4. Lege dieselbe Tabelle als "

# --- 3. Update the cached date placeholder text -----------------------------
$today = "12/17/2024"

$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $mshp = $m.Shapes.Item($i)
    if ($mshp.PlaceholderFormat.Type -eq 16) {
        $mtr = $mshp.TextFrame.TextRange
        $mlen = $mtr.Text.Length
        $msub = $mtr.Characters(1, $mlen)
        $msub.Text = $today
    }
}

$layouts = $m.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $cl = $layouts.Item($j)
    for ($i = 1; $i -le $cl.Shapes.Count; $i++) {
        $lshp = $cl.Shapes.Item($i)
        if ($lshp.PlaceholderFormat.Type -eq 16) {
            $ltr = $lshp.TextFrame.TextRange
            $llen = $ltr.Text.Length
            $lsub = $ltr.Characters(1, $llen)
            $lsub.Text = $today
        }
    }
}

Write-Output "done"
